# Reorder the comma-separated "Recorded By" names in column G so that the
# first author in the list is moved to the end (rotate left by one).
# Example: "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("G$row")
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -gt 1) {
        $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
        $newVal = $rotated -join ", "
        $cell.Value2 = $newVal
    }
}
